# Automatische test-sync: 2025-06-26 23:17:50
#
# Adds a new "Testmail #4" row (row 36) to the Logs sheet, widens the
# conditional-formatting ranges that tracked the old last row (35) to
# include the new last row (36), and updates the Dashboard category
# counts/order to reflect the new "Retour / Terugbetaling" entry.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Logs sheet: append the new test-mail row
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A36").Value = "Hoe kan ik iets retourneren?"
$logs.Range("B36").Value = "mailmind.test@zohomail.eu"
$logs.Range("C36").Value = "Testmail #4: Hoe kan ik iets retourneren?"
$logs.Range("D36").Value = "Retour / Terugbetaling"
$logs.Range("E36").Value = "Beste klant,`nBedankt voor je bericht. Als je een artikel wilt retourneren, kun je dit doen door contact op te nemen met onze klantenservice via support@bedrijf.nl. Zij zullen je verder begeleiden bij het retourproces en eventuele vragen beantwoorden.`nBedankt voor je begrip en medewerking.`nMet vriendelijke groet,`n[Bedrijfsnaam]"
$logs.Range("F36").Value = "2025-06-26 23:17:40"
$logs.Range("G36").Value = "Ja"
$logs.Range("H36").Value = "Nee"
$logs.Range("I36").Value = "Ja"

# ---------------------------------------------------------------------
# 2) Logs sheet: extend the conditional-formatting ranges so the new
#    row is covered (D/G/H/I 2:35 -> 2:36)
# ---------------------------------------------------------------------
$logs.Range("D2:D35").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D36"))
$logs.Range("G2:G35").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G36"))
$logs.Range("H2:H35").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H36"))
$logs.Range("I2:I35").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I36"))

# ---------------------------------------------------------------------
# 3) Dashboard sheet: "Retour / Terugbetaling" now ties "Productinformatie"
#    at 3, and swaps ahead of it in the ranking
# ---------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Range("A4").Value = "Retour / Terugbetaling"
$dashboard.Range("B4").Value = 3

$dashboard.Range("A5").Value = "Productinformatie"
$dashboard.Range("B5").Value = 3
